# Updated cryptos list on Wed Mar 15 09:10:21 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold numeric-looking text (e.g. "24.905.55",
# "  +2.29%  ") that must stay text, not be reinterpreted as numbers.
# Mark the whole data range as Text before writing any values so Excel
# keeps every cell (touched or not) as a text value.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Simple price / volume updates ---
$ws.Range("D2").Value = "24.905.55"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "1.706.85"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "312.83"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.3745"
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("D8").Value = "49.33"
$ws.Range("E8").Value = "  +4.05%  "
$ws.Range("D9").Value = "0.3440"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +5.86%  "
$ws.Range("D11").Value = "0.07523"
$ws.Range("E11").Value = "  +4.39%  "
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "21.18"
$ws.Range("E13").Value = "  +5.89%  "
$ws.Range("D14").Value = "6.363"
$ws.Range("E14").Value = "  +4.05%  "
$ws.Range("D15").Value = "7.050"
$ws.Range("E15").Value = "  +5.29%  "
$ws.Range("D16").Value = "1.708.62"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "0.06732"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D20").Value = "84.03"
$ws.Range("E20").Value = "  +4.52%  "
$ws.Range("D21").Value = "17.38"
$ws.Range("E21").Value = "  +6.19%  "
$ws.Range("D22").Value = "6.382"
$ws.Range("E22").Value = "  +4.96%  "
$ws.Range("D23").Value = "13.23"
$ws.Range("E23").Value = "  +8.28%  "
$ws.Range("D24").Value = "24.897.18"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("D25").Value = "2.449"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("D26").Value = "2.800"
$ws.Range("E26").Value = "  +6.24%  "
$ws.Range("D27").Value = "20.40"
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").Value = "132.75"
$ws.Range("E29").Value = "  +3.97%  "

# --- Rows 30/31 swap: ImmutableX now ranks above WrappedliquidstakedEther2.0 ---
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "1.266"
$ws.Range("E30").Value = "  +31.30%  "

$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "1.896.40"
$ws.Range("E31").Value = "  +2.00%  "

$ws.Range("D32").Value = "6.820"
$ws.Range("E32").Value = "  +9.23%  "
$ws.Range("D33").Value = "4.231"
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("D34").Value = "13.80"
$ws.Range("E34").Value = "  +13.04%  "

# --- Rows 35/36 swap: Stellar now ranks above WEMIXTOKEN ---
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.08781"
$ws.Range("E35").Value = "  +4.23%  "

$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "1.774"
$ws.Range("E36").Value = "  +5.87%  "

$ws.Range("D37").Value = "5.630"
$ws.Range("E37").Value = "  +6.65%  "
$ws.Range("D38").Value = "0.06652"
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("D39").Value = "9.174"
$ws.Range("E39").Value = "  +6.46%  "
$ws.Range("D40").Value = "0.02421"
$ws.Range("E40").Value = "  +4.96%  "
$ws.Range("D41").Value = "0.2261"
$ws.Range("E41").Value = "  +8.86%  "
$ws.Range("D42").Value = "1.272"
$ws.Range("E42").Value = "  +2.58%  "
$ws.Range("D43").Value = "0.6460"
$ws.Range("E43").Value = "  +6.68%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "13.92"
$ws.Range("E45").Value = "  +8.06%  "
$ws.Range("D46").Value = "0.6163"
$ws.Range("E46").Value = "  +5.33%  "
$ws.Range("D47").Value = "3.838"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").Value = "2.118"
$ws.Range("E48").Value = "  +5.55%  "
$ws.Range("D49").Value = "129.99"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("D50").Value = "0.07319"
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").Value = "80.26"
$ws.Range("E51").Value = "  +6.32%  "
